$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(3).Insert()

$ws.Cells.Item(3, 1).Value = 9
$ws.Cells.Item(3, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(3, 3).Value = "Metropolitana"
$ws.Cells.Item(3, 4).Value = 44550
$ws.Cells.Item(3, 5).Value = 13
$ws.Cells.Item(3, 6).Value = "Fruta"
$ws.Cells.Item(3, 7).Value = 100101
$ws.Cells.Item(3, 8).Value = "Berries"
$ws.Cells.Item(3, 9).Value = 100101006
$ws.Cells.Item(3, 10).Value = "Breva"
$ws.Cells.Item(3, 11).Value = "Sin especificar"
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 60
$ws.Cells.Item(3, 14).Value = 24000
$ws.Cells.Item(3, 15).Value = 24000
$ws.Cells.Item(3, 16).Value = 24000
$ws.Cells.Item(3, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(3, 18).Value = "Región Metropolitana"
$ws.Cells.Item(3, 19).Value = 3429
$ws.Cells.Item(3, 20).Value = 7
